$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header in column T (20th column) on row 1, matching the
# style of the existing header cells (copy style/number format from S1).
$ws.Range("T1").Value = "Nilai Pembayaran Zakat"
$ws.Range("T1").NumberFormat = $ws.Range("S1").NumberFormat

# Update the active selection to match the authored workbook state.
$ws.Range("R7").Select()
